$d = $word.ActiveDocument

# 1. Version 0.1 -> 1.2.5
$d.Content.Find.Execute("0.1", $true, $false, $false, $false, $false, $true, 1, $false, "1.2.5", 2) | Out-Null

# 2. Creation -> Update
$d.Content.Find.Execute("Creation", $true, $false, $false, $false, $false, $true, 1, $false, "Update", 2) | Out-Null

# 3. Date 17/02/2023 -> 31/05/2023
$d.Content.Find.Execute("17/02/2023", $true, $false, $false, $false, $false, $true, 1, $false, "31/05/2023", 2) | Out-Null

# 4. "O usuario devidamente autenticado e na tela de listagem de empenhos" -> with accent + period
$d.Content.Find.Execute("O usuario devidamente autenticado e na tela de listagem de empenhos", $true, $false, $false, $false, $false, $true, 1, $false, "O usuário devidamente autenticado e na tela de listagem de empenhos.", 2) | Out-Null

# 5. Step 2 text about listagem
$d.Content.Find.Execute("2. System Exibe a lista de solicitações aguardando serem empenhadas ordenado pelo numero da diaria em ordem crescente. ", $true, $false, $false, $false, $false, $true, 1, $false, "2. System Exibe a lista de solicitações aguardando serem empenhadas, de todos os servidores, ordenado pelo número da diária em ordem crescente. ", 2) | Out-Null

# 6. Apresenta a tela de Detalhar Diárias (add period)
$d.Content.Find.Execute("Apresenta a tela de Detalhar Diárias ", $true, $false, $false, $false, $false, $true, 1, $false, "Apresenta a tela de Detalhar Diárias. ", 2) | Out-Null

# 7. filtra -> Filtra + period
$d.Content.Find.Execute("1. Chefe/Beneficiário filtra a listagem por registros cujos beneficiários não possuem número do credor ", $true, $false, $false, $false, $false, $true, 1, $false, "1. Chefe/Beneficiário Filtra a listagem por registros cujos beneficiários não possuem número do credor. ", 2) | Out-Null

# 8. Registrar Empenho (add period)
$d.Content.Find.Execute("2. System Apresenta a tela de Registrar Empenho ", $true, $false, $false, $false, $false, $true, 1, $false, "2. System Apresenta a tela de Registrar Empenho. ", 2) | Out-Null
